$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells whose new values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts them to the Number type
# (these are meant to stay plain text, matching the rest of the column).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '66.705.12'
$ws.Range("E2").Value = '  -3.79%  '
$ws.Range("D3").Value = '3.318.09'
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '573.68'
$ws.Range("E5").Value = '  -2.97%  '
$ws.Range("D6").Value = '182.88'
$ws.Range("E6").Value = '  -5.16%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '0.602'
$ws.Range("E8").Value = '  -0.81%  '
$ws.Range("D9").Value = '0.130'
$ws.Range("E9").Value = '  -2.94%  '
$ws.Range("E10").Value = '  -1.91%  '
$ws.Range("E11").Value = '  -4.62%  '
$ws.Range("D12").Value = '3.893.96'
$ws.Range("E12").Value = '  -0.91%  '
$ws.Range("E13").Value = '  -0.54%  '
$ws.Range("D14").Value = '27.23'
$ws.Range("E14").Value = '  -3.31%  '
$ws.Range("D15").Value = '66.703.11'
$ws.Range("E15").Value = '  -3.78%  '
$ws.Range("E16").Value = '  -2.27%  '
$ws.Range("D17").Value = '3.289.98'
$ws.Range("E17").Value = '  -0.58%  '
$ws.Range("D18").Value = '13.78'
$ws.Range("E18").Value = '  +0.56%  '
$ws.Range("D19").Value = '437.36'
$ws.Range("E19").Value = '  +2.62%  '
$ws.Range("E20").Value = '  -2.32%  '
$ws.Range("E21").Value = '  -0.98%  '
$ws.Range("D22").Value = '73.89'
$ws.Range("E22").Value = '  +0.82%  '
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("E25").Value = '  -1.94%  '
$ws.Range("E26").Value = '  +1.15%  '
$ws.Range("D27").Value = '9.08'
$ws.Range("E27").Value = '  -5.19%  '
$ws.Range("E28").Value = '  -0.60%  '
$ws.Range("E29").Value = '  -1.68%  '
$ws.Range("D30").Value = '22.92'
$ws.Range("E30").Value = '  -0.43%  '
$ws.Range("D31").Value = '5.34'
$ws.Range("E31").Value = '  -4.28%  '
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("D33").Value = '6.78'
$ws.Range("E33").Value = '  -2.84%  '
$ws.Range("E34").Value = '  -3.57%  '
$ws.Range("E35").Value = '  -0.94%  '
$ws.Range("D36").Value = '160.24'
$ws.Range("E36").Value = '  -2.75%  '
$ws.Range("D37").Value = '27.42'
$ws.Range("E37").Value = '  +1.43%  '
$ws.Range("D38").Value = '1.86'
$ws.Range("E38").Value = '  -3.18%  '
$ws.Range("D39").Value = '2.810.64'
$ws.Range("E39").Value = '  +2.13%  '
$ws.Range("D40").Value = '0.791'
$ws.Range("E40").Value = '  -2.22%  '
$ws.Range("E41").Value = '  -2.24%  '
$ws.Range("D42").Value = '6.24'
$ws.Range("E42").Value = '  -3.75%  '
$ws.Range("E43").Value = '  -1.34%  '
$ws.Range("E44").Value = '  -2.37%  '
$ws.Range("D45").Value = '24.32'
$ws.Range("E45").Value = '  -4.05%  '
$ws.Range("E46").Value = '  -6.38%  '
$ws.Range("D47").Value = '318.85'
$ws.Range("E47").Value = '  -7.97%  '
$ws.Range("E48").Value = '  -2.82%  '
$ws.Range("D49").Value = '0.984'
$ws.Range("E49").Value = '  -2.28%  '
$ws.Range("D50").Value = '6.19'
$ws.Range("E50").Value = '  -1.40%  '
$ws.Range("D51").Value = '0.0998'
$ws.Range("E51").Value = '  -1.53%  '
